# For each 4-row year block (A,B,C,D rows), the B-row and C-row need to be
# swapped (e.g. "2000年B" / "2000年C" trade places, and so on for every year
# from 2000 through 2019). Columns F and G (the duplicate/derived
# "产销率" and "销售量" quarter columns) are then removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Staging row, well outside the A1:G81 data range, used as scratch space for
# the 3-step swap (cut row1 -> staging, cut row2 -> row1, cut staging -> row2).
$stageRow = 1000
$stageRange = "A" + $stageRow + ":E" + $stageRow

$firstDataRow = 2
$lastDataRow = 81
$blockSize = 4

for ($base = $firstDataRow; $base -le $lastDataRow; $base += $blockSize) {
    $rowB = $base + 1
    $rowC = $base + 2

    $rangeB = "A" + $rowB + ":E" + $rowB
    $rangeC = "A" + $rowC + ":E" + $rowC

    $ws.Range($rangeB).Cut($ws.Range($stageRange)) | Out-Null
    $ws.Range($rangeC).Cut($ws.Range($rangeB)) | Out-Null
    $ws.Range($stageRange).Cut($ws.Range($rangeC)) | Out-Null
}

# Remove the scratch row entirely so no stray formatting/dimension residue
# from the staging cells is left behind.
$ws.Rows($stageRow).Delete() | Out-Null

# Drop the now-redundant F (产销率) and G (销售量) quarter-value columns,
# shifting nothing needs to move in since they are the last two columns.
$ws.Range("F1:G81").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft) | Out-Null
